$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 8.284530645224702
$ws.Cells.Item(2, 3).Value = 4.18747655020839
$ws.Cells.Item(2, 5).Value = 20.94896542863705
$ws.Cells.Item(2, 6).Value = 39.9813710967377
$ws.Cells.Item(2, 7).Value = 3.641128967802692
$ws.Cells.Item(2, 9).Value = 21.1248927433636
$ws.Cells.Item(2, 10).Value = 8.156773086012334
$ws.Cells.Item(2, 11).Value = 8.318792822958827
$ws.Cells.Item(2, 13).Value = 17.8599260892343
$ws.Cells.Item(2, 14).Value = 18.92386595907361
$ws.Cells.Item(2, 15).Value = 22.466435791962
$ws.Cells.Item(3, 2).Value = 8.020266343369824
$ws.Cells.Item(3, 3).Value = 4.012130654381945
$ws.Cells.Item(3, 5).Value = 20.84374269683949
$ws.Cells.Item(3, 6).Value = 39.9234818999803
$ws.Cells.Item(3, 7).Value = 3.642654102948329
$ws.Cells.Item(3, 9).Value = 21.21141868567004
$ws.Cells.Item(3, 10).Value = 8.16891612235993
$ws.Cells.Item(3, 11).Value = 8.152538795431846
$ws.Cells.Item(3, 13).Value = 17.74538491153528
$ws.Cells.Item(3, 14).Value = 18.98054857589916
$ws.Cells.Item(3, 15).Value = 22.54634562468903
$ws.Cells.Item(4, 2).Value = 7.854691225476741
$ws.Cells.Item(4, 3).Value = 3.899596870400606
$ws.Cells.Item(4, 5).Value = 20.78329407065045
$ws.Cells.Item(4, 6).Value = 39.89767810589929
$ws.Cells.Item(4, 7).Value = 3.643640166525413
$ws.Cells.Item(4, 9).Value = 21.26814182455064
$ws.Cells.Item(4, 10).Value = 8.176757767591377
$ws.Cells.Item(4, 11).Value = 8.049927299828383
$ws.Cells.Item(4, 13).Value = 17.67787037404031
$ws.Cells.Item(4, 14).Value = 19.01700449518246
$ws.Cells.Item(4, 15).Value = 22.59970267391553
$ws.Cells.Item(5, 2).Value = 7.786490245680477
$ws.Cells.Item(5, 3).Value = 3.85255093406799
$ws.Cells.Item(5, 5).Value = 20.75972551920388
$ws.Cells.Item(5, 6).Value = 39.8896199065052
$ws.Cells.Item(5, 7).Value = 3.644054512068789
$ws.Cells.Item(5, 9).Value = 21.29216125722008
$ws.Cells.Item(5, 10).Value = 8.180050578712391
$ws.Cells.Item(5, 11).Value = 8.008035116901741
$ws.Cells.Item(5, 13).Value = 17.65108819213587
$ws.Cells.Item(5, 14).Value = 19.03227733091112
$ws.Cells.Item(5, 15).Value = 22.62252406627532
$ws.Cells.Item(6, 2).Value = 7.77512482896433
$ws.Cells.Item(6, 3).Value = 3.844668467004004
$ws.Cells.Item(6, 5).Value = 20.75587681474172
$ws.Cells.Item(6, 6).Value = 39.888430433808
$ws.Cells.Item(6, 7).Value = 3.644124070927485
$ws.Cells.Item(6, 9).Value = 21.29620428588305
$ws.Cells.Item(6, 10).Value = 8.180603230773741
$ws.Cells.Item(6, 11).Value = 8.001075965567013
$ws.Cells.Item(6, 13).Value = 17.64668580050056
$ws.Cells.Item(6, 14).Value = 19.03483857988144
$ws.Cells.Item(6, 15).Value = 22.62637860849038
$ws.Cells.Item(7, 2).Value = 7.853774243797648
$ws.Cells.Item(7, 3).Value = 3.89896714735302
$ws.Cells.Item(7, 5).Value = 20.78297188174345
$ws.Cells.Item(7, 6).Value = 39.89755947304224
$ws.Cells.Item(7, 7).Value = 3.643645703804477
$ws.Cells.Item(7, 9).Value = 21.26846209744385
$ws.Cells.Item(7, 10).Value = 8.176801781361764
$ws.Cells.Item(7, 11).Value = 8.049362564154237
$ws.Cells.Item(7, 13).Value = 17.67750619286143
$ws.Cells.Item(7, 14).Value = 19.01720878097822
$ws.Cells.Item(7, 15).Value = 22.60000608833329
$ws.Cells.Item(8, 2).Value = 8.194161520489306
$ws.Cells.Item(8, 3).Value = 4.128048903200201
$ws.Cells.Item(8, 5).Value = 20.91183340234831
$ws.Cells.Item(8, 6).Value = 39.9593945554733
$ws.Cells.Item(8, 7).Value = 3.641644557864814
$ws.Cells.Item(8, 9).Value = 21.15398061547436
$ws.Cells.Item(8, 10).Value = 8.16088010199854
$ws.Cells.Item(8, 11).Value = 8.261611921522205
$ws.Cells.Item(8, 13).Value = 17.81986226652459
$ws.Cells.Item(8, 14).Value = 18.94306786229877
$ws.Cells.Item(8, 15).Value = 22.49309709307482
$ws.Cells.Item(9, 2).Value = 8.831119402212916
$ws.Cells.Item(9, 3).Value = 4.53725552807279
$ws.Cells.Item(9, 5).Value = 21.19662217080783
$ws.Cells.Item(9, 6).Value = 40.15750692602039
$ws.Cells.Item(9, 7).Value = 3.638112320169435
$ws.Cells.Item(9, 9).Value = 20.95801041163014
$ws.Cells.Item(9, 10).Value = 8.132706170357288
$ws.Cells.Item(9, 11).Value = 8.671231355985059
$ws.Cells.Item(9, 13).Value = 18.12029654821953
$ws.Cells.Item(9, 14).Value = 18.81073440319502
$ws.Cells.Item(9, 15).Value = 22.31756306049556
$ws.Cells.Item(10, 2).Value = 9.275403844641355
$ws.Cells.Item(10, 3).Value = 4.81198219750892
$ws.Cells.Item(10, 5).Value = 21.42411762455432
$ws.Cells.Item(10, 6).Value = 40.34920139870955
$ws.Cells.Item(10, 7).Value = 3.635753715492207
$ws.Cells.Item(10, 9).Value = 20.83141684524853
$ws.Cells.Item(10, 10).Value = 8.113847558230697
$ws.Cells.Item(10, 11).Value = 8.965101384268278
$ws.Cells.Item(10, 13).Value = 18.35251833280422
$ws.Cells.Item(10, 14).Value = 18.72138955247878
$ws.Cells.Item(10, 15).Value = 22.20947407190578
$ws.Cells.Item(11, 2).Value = 9.471470656945765
$ws.Cells.Item(11, 3).Value = 4.931092723346225
$ws.Cells.Item(11, 5).Value = 21.53127327180644
$ws.Cells.Item(11, 6).Value = 40.44623177401988
$ws.Cells.Item(11, 7).Value = 3.634731565398033
$ws.Cells.Item(11, 9).Value = 20.77760084846303
$ws.Cells.Item(11, 10).Value = 8.105664238359072
$ws.Cells.Item(11, 11).Value = 9.096655754429843
$ws.Cells.Item(11, 13).Value = 18.46032690413069
$ws.Cells.Item(11, 14).Value = 18.68243859226844
$ws.Cells.Item(11, 15).Value = 22.16485122040796
$ws.Cells.Item(12, 2).Value = 9.544782517843869
$ws.Cells.Item(12, 3).Value = 4.97533756375404
$ws.Cells.Item(12, 5).Value = 21.57235043217558
$ws.Cells.Item(12, 6).Value = 40.48436792230598
$ws.Cells.Item(12, 7).Value = 3.634351768177361
$ws.Cells.Item(12, 9).Value = 20.75776468818457
$ws.Cells.Item(12, 10).Value = 8.102622028164243
$ws.Cells.Item(12, 11).Value = 9.146119712048376
$ws.Cells.Item(12, 13).Value = 18.50143447183877
$ws.Cells.Item(12, 14).Value = 18.66793100915534
$ws.Cells.Item(12, 15).Value = 22.14860900426671
$ws.Cells.Item(13, 2).Value = 9.529036023740748
$ws.Cells.Item(13, 3).Value = 4.965847115611814
$ws.Cells.Item(13, 5).Value = 21.56348195852066
$ws.Cells.Item(13, 6).Value = 40.47609301845772
$ws.Cells.Item(13, 7).Value = 3.63443324157249
$ws.Cells.Item(13, 9).Value = 20.76201261317522
$ws.Cells.Item(13, 10).Value = 8.103274707568287
$ws.Cells.Item(13, 11).Value = 9.135483169055222
$ws.Cells.Item(13, 13).Value = 18.49256911311218
$ws.Cells.Item(13, 14).Value = 18.67104471735882
$ws.Cells.Item(13, 15).Value = 22.15207787818225
$ws.Cells.Item(14, 2).Value = 9.477521145132723
$ws.Cells.Item(14, 3).Value = 4.934750087851127
$ws.Cells.Item(14, 5).Value = 21.5346428624205
$ws.Cells.Item(14, 6).Value = 40.44934146681291
$ws.Cells.Item(14, 7).Value = 3.634700173772136
$ws.Cells.Item(14, 9).Value = 20.77595803504048
$ws.Cells.Item(14, 10).Value = 8.105412820233843
$ws.Cells.Item(14, 11).Value = 9.100732512543397
$ws.Cells.Item(14, 13).Value = 18.46370334938574
$ws.Cells.Item(14, 14).Value = 18.68124019472672
$ws.Cells.Item(14, 15).Value = 22.16350181757183
$ws.Cells.Item(15, 2).Value = 9.445843282803953
$ws.Cells.Item(15, 3).Value = 4.915589856410815
$ws.Cells.Item(15, 5).Value = 21.51704231589276
$ws.Cells.Item(15, 6).Value = 40.43313612088236
$ws.Cells.Item(15, 7).Value = 3.634864623040424
$ws.Cells.Item(15, 9).Value = 20.78457070752692
$ws.Cells.Item(15, 10).Value = 8.106729844341976
$ws.Cells.Item(15, 11).Value = 9.079399445992758
$ws.Cells.Item(15, 13).Value = 18.44605819114044
$ws.Cells.Item(15, 14).Value = 18.68751674197203
$ws.Cells.Item(15, 15).Value = 22.17058472081809
$ws.Cells.Item(16, 2).Value = 9.262463076779587
$ws.Cells.Item(16, 3).Value = 4.804078570792884
$ws.Cells.Item(16, 5).Value = 21.41718619476991
$ws.Cells.Item(16, 6).Value = 40.34305629922278
$ws.Cells.Item(16, 7).Value = 3.635821534463816
$ws.Cells.Item(16, 9).Value = 20.83500978162999
$ws.Cells.Item(16, 10).Value = 8.114390296044386
$ws.Cells.Item(16, 11).Value = 8.956457102653676
$ws.Cells.Item(16, 13).Value = 18.34551396291776
$ws.Cells.Item(16, 14).Value = 18.72396905392609
$ws.Cells.Item(16, 15).Value = 22.21248192091668
$ws.Cells.Item(17, 2).Value = 9.148368998640837
$ws.Cells.Item(17, 3).Value = 4.734155444183583
$ws.Cells.Item(17, 5).Value = 21.3568474879165
$ws.Cells.Item(17, 6).Value = 40.29029850388826
$ws.Cells.Item(17, 7).Value = 3.636421552534849
$ws.Cells.Item(17, 9).Value = 20.86691888721554
$ws.Cells.Item(17, 10).Value = 8.119190871622825
$ws.Cells.Item(17, 11).Value = 8.880456603849076
$ws.Cells.Item(17, 13).Value = 18.28436841607882
$ws.Cells.Item(17, 14).Value = 18.74676407880648
$ws.Cells.Item(17, 15).Value = 22.2393503226018
$ws.Cells.Item(18, 2).Value = 9.082179996455372
$ws.Cells.Item(18, 3).Value = 4.693386353579532
$ws.Cells.Item(18, 5).Value = 21.32248911522993
$ws.Cells.Item(18, 6).Value = 40.26087998468271
$ws.Cells.Item(18, 7).Value = 3.636771449812615
$ws.Cells.Item(18, 9).Value = 20.88562717954098
$ws.Cells.Item(18, 10).Value = 8.121989281983737
$ws.Cells.Item(18, 11).Value = 8.836545389154548
$ws.Cells.Item(18, 13).Value = 18.24940511784273
$ws.Cells.Item(18, 14).Value = 18.76003454002542
$ws.Cells.Item(18, 15).Value = 22.25523224383347
$ws.Cells.Item(19, 2).Value = 9.059674606059165
$ws.Cells.Item(19, 3).Value = 4.679488527605315
$ws.Cells.Item(19, 5).Value = 21.31091636054783
$ws.Cells.Item(19, 6).Value = 40.25107907087513
$ws.Cells.Item(19, 7).Value = 3.636890741528689
$ws.Cells.Item(19, 9).Value = 20.89202244990116
$ws.Cells.Item(19, 10).Value = 8.122943180236934
$ws.Cells.Item(19, 11).Value = 8.821645316282959
$ws.Cells.Item(19, 13).Value = 18.23760340966759
$ws.Cells.Item(19, 14).Value = 18.7645550943213
$ws.Cells.Item(19, 15).Value = 22.26068303769548
$ws.Cells.Item(20, 2).Value = 9.160573494287219
$ws.Cells.Item(20, 3).Value = 4.741656051535542
$ws.Cells.Item(20, 5).Value = 21.36323494535721
$ws.Cells.Item(20, 6).Value = 40.2958189264108
$ws.Cells.Item(20, 7).Value = 3.636357184842849
$ws.Cells.Item(20, 9).Value = 20.86348536274049
$ws.Cells.Item(20, 10).Value = 8.118675989091688
$ws.Cells.Item(20, 11).Value = 8.88856779960158
$ws.Cells.Item(20, 13).Value = 18.29085635355798
$ws.Cells.Item(20, 14).Value = 18.74432102353529
$ws.Cells.Item(20, 15).Value = 22.23644583493879
$ws.Cells.Item(21, 2).Value = 9.492678160655259
$ws.Cells.Item(21, 3).Value = 4.94390748918451
$ws.Cells.Item(21, 5).Value = 21.54310026855943
$ws.Cells.Item(21, 6).Value = 40.45716141659732
$ws.Cells.Item(21, 7).Value = 3.634621572317282
$ws.Cells.Item(21, 9).Value = 20.77184719228865
$ws.Cells.Item(21, 10).Value = 8.104783269589788
$ws.Cells.Item(21, 11).Value = 9.110949556410617
$ws.Cells.Item(21, 13).Value = 18.47217448345843
$ws.Cells.Item(21, 14).Value = 18.67823896842302
$ws.Cells.Item(21, 15).Value = 22.16012852739816
$ws.Cells.Item(22, 2).Value = 9.704255990280835
$ws.Cells.Item(22, 3).Value = 5.071073568536592
$ws.Cells.Item(22, 5).Value = 21.66355153391635
$ws.Cells.Item(22, 6).Value = 40.57071581288413
$ws.Cells.Item(22, 7).Value = 3.633529601578042
$ws.Cells.Item(22, 9).Value = 20.7151203427028
$ws.Cells.Item(22, 10).Value = 8.096033576656485
$ws.Cells.Item(22, 11).Value = 9.254212728650964
$ws.Cells.Item(22, 13).Value = 18.59231094800872
$ws.Cells.Item(22, 14).Value = 18.63646234618552
$ws.Cells.Item(22, 15).Value = 22.11407216059766
$ws.Cells.Item(23, 2).Value = 9.591853452062351
$ws.Cells.Item(23, 3).Value = 5.00366643279834
$ws.Cells.Item(23, 5).Value = 21.59900853381454
$ws.Cells.Item(23, 6).Value = 40.50937507527047
$ws.Cells.Item(23, 7).Value = 3.634108543104156
$ws.Cells.Item(23, 9).Value = 20.7451068919606
$ws.Cells.Item(23, 10).Value = 8.100673337686557
$ws.Cells.Item(23, 11).Value = 9.17795513260443
$ws.Cells.Item(23, 13).Value = 18.52805203520315
$ws.Cells.Item(23, 14).Value = 18.65863048573331
$ws.Cells.Item(23, 15).Value = 22.13830313656907
$ws.Cells.Item(24, 2).Value = 9.155057689386469
$ws.Cells.Item(24, 3).Value = 4.73826679910912
$ws.Cells.Item(24, 5).Value = 21.3603461422296
$ws.Cells.Item(24, 6).Value = 40.29332029903351
$ws.Cells.Item(24, 7).Value = 3.636386270089316
$ws.Cells.Item(24, 9).Value = 20.86503652744483
$ws.Cells.Item(24, 10).Value = 8.118908647554475
$ws.Cells.Item(24, 11).Value = 8.884901401920544
$ws.Cells.Item(24, 13).Value = 18.28792256267582
$ws.Cells.Item(24, 14).Value = 18.74542501383764
$ws.Cells.Item(24, 15).Value = 22.23775759917044
$ws.Cells.Item(25, 2).Value = 8.662631101649849
$ws.Cells.Item(25, 3).Value = 4.43101178907412
$ws.Cells.Item(25, 5).Value = 21.11627299543317
$ws.Cells.Item(25, 6).Value = 40.09574941207794
$ws.Cells.Item(25, 7).Value = 3.639026171260545
$ws.Cells.Item(25, 9).Value = 21.00797216454887
$ws.Cells.Item(25, 10).Value = 8.140003439838564
$ws.Cells.Item(25, 11).Value = 8.561449483329469
$ws.Cells.Item(25, 13).Value = 18.03689262841466
$ws.Cells.Item(25, 14).Value = 18.84514452930096
$ws.Cells.Item(25, 15).Value = 22.36138927948754
